# Insert a new weekly record at row 32 (pushing the existing rows 32-52 down
# to 33-53) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32..52 down by one row to make room for the new record.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new weekly entry.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44729
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112037
$ws.Range("G32").Value = "Cebollín"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 220
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 6500
$ws.Range("M32").Value = 6273
$ws.Range("N32").Value = "$/paquete 36 unidades"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 174
$ws.Range("Q32").Value = 36
$ws.Range("R32").Value = "Hortaliza"
